# Update cryptocurrency price/volume data per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.219.75"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'3.673.38"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'674.10"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'157.62"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "'0.146"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").Value = "  -5.67%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "'0.0000231"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("D13").Value = "'4.290.61"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "'32.16"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "'3.670.06"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "'69.173.43"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'16.00"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").Value = "'466.91"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'3.818.77"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D26").Value = "'10.87"
$ws.Range("E26").Value = "  -5.50%  "
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = "  -7.99%  "
$ws.Range("D28").Value = "'8.97"
$ws.Range("E28").Value = "  -6.50%  "
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  -6.29%  "
$ws.Range("E31").Value = "  -4.08%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'26.83"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  -5.47%  "
$ws.Range("D35").Value = "'3.663.26"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'0.159"
$ws.Range("E36").Value = "  -5.15%  "
$ws.Range("D37").Value = "'8.11"
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("D38").Value = "'6.14"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("D43").Value = "'171.99"
$ws.Range("E43").Value = "  +7.74%  "
$ws.Range("D44").Value = "'0.940"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").Value = "'47.53"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "'0.000275"
$ws.Range("E46").Value = "  -5.05%  "
$ws.Range("D47").Value = "'2.66"
$ws.Range("E47").Value = "  -6.67%  "
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("D49").Value = "'27.17"
$ws.Range("E49").Value = "  -9.64%  "
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "'7.75"
$ws.Range("E51").Value = "  -3.62%  "
